$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2383.353
$ws.Range("I28").Value = 2975
$ws.Range("J28").Value = 460.5
$ws.Range("K28").Value = 2975
$ws.Range("L28").Value = 460.5
$ws.Range("M28").Value = -2490
$ws.Range("N28").Value = -1430.5
$ws.Range("H40").Value = 6183.4165
$ws.Range("I40").Value = 3800.5
$ws.Range("J40").Value = 7374.875
$ws.Range("K40").Value = 3800.5
$ws.Range("L40").Value = 7374.875
$ws.Range("M40").Value = -3625.5
$ws.Range("N40").Value = -7724.875
$ws.Range("H41").Value = 329.66666
$ws.Range("I41").Value = 66.333336
$ws.Range("J41").Value = 593
$ws.Range("K41").Value = 66.333336
$ws.Range("L41").Value = 593
$ws.Range("M41").Value = 373.666664
$ws.Range("N41").Value = -1473
$ws.Range("H43").Value = 67812.336
$ws.Range("I43").Value = 3000
$ws.Range("K43").Value = 3000
$ws.Range("M43").Value = -2931
$ws.Range("H53").Value = 309.63635
$ws.Range("I53").Value = 380.92856
$ws.Range("J53").Value = 184.875
$ws.Range("K53").Value = 380.92856
$ws.Range("L53").Value = 184.875
$ws.Range("M53").Value = 256.07144
$ws.Range("N53").Value = -1458.875
$ws.Range("H64").Value = 5000
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5496
$ws.Range("H67").Value = 5000
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6716
$ws.Range("H69").Value = 8142.7144
$ws.Range("I69").Value = 7110.8887
$ws.Range("K69").Value = 21332.6661
$ws.Range("M69").Value = -20458.6661
$ws.Range("H72").Value = 8142.7144
$ws.Range("I72").Value = 7110.8887
$ws.Range("K72").Value = 63997.99830000001
$ws.Range("M72").Value = -59629.99830000001
$ws.Range("H76").Value = 1998
$ws.Range("J76").Value = 1998
$ws.Range("L76").Value = 1998
$ws.Range("N76").Value = -2628
$ws.Range("H79").Value = 1998
$ws.Range("J79").Value = 1998
$ws.Range("L79").Value = 1998
$ws.Range("N79").Value = -4182
$ws.Range("H97").Value = 10982
$ws.Range("J97").Value = 11178.6
$ws.Range("L97").Value = 33535.8
$ws.Range("N97").Value = -34527.8
$ws.Range("H113").Value = 8179.143
$ws.Range("I113").Value = 6621.2856
$ws.Range("J113").Value = 9737
$ws.Range("K113").Value = 6621.2856
$ws.Range("L113").Value = 9737
$ws.Range("M113").Value = -3367.2856
$ws.Range("N113").Value = -16245
$ws.Range("H135").Value = 698.1818
$ws.Range("I135").Value = 631.3333
$ws.Range("K135").Value = 5681.9997
$ws.Range("M135").Value = -3146.9997
$ws.Range("H137").Value = 1766.7858
$ws.Range("I137").Value = 1099.75
$ws.Range("J137").Value = 2033.6
$ws.Range("K137").Value = 3299.25
$ws.Range("L137").Value = 6100.799999999999
$ws.Range("M137").Value = -749.25
$ws.Range("N137").Value = -11200.8
$ws.Range("H138").Value = 5281.2393
$ws.Range("I138").Value = 1480.3334
$ws.Range("J138").Value = 5546.4185
$ws.Range("K138").Value = 4441.0002
$ws.Range("L138").Value = 16639.2555
$ws.Range("M138").Value = 698.9997999999996
$ws.Range("N138").Value = -26919.2555
$ws.Range("H141").Value = 5366.25
$ws.Range("I141").Value = 3155
$ws.Range("J141").Value = 12000
$ws.Range("K141").Value = 9465
$ws.Range("L141").Value = 36000
$ws.Range("M141").Value = -4285
$ws.Range("N141").Value = -46360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5598.4517
$ws.Range("I32").Value = 3541.5833
$ws.Range("K32").Value = 3541.5833
$ws.Range("M32").Value = -3254.5833
$ws.Range("H41").Value = 17338.572
$ws.Range("I41").Value = 12040.333
$ws.Range("J41").Value = 21312.25
$ws.Range("K41").Value = 12040.333
$ws.Range("L41").Value = 21312.25
$ws.Range("M41").Value = -11626.333
$ws.Range("N41").Value = -22140.25
$ws.Range("H74").Value = 1829.8182
$ws.Range("I74").Value = 1792.1111
$ws.Range("K74").Value = 1792.1111
$ws.Range("M74").Value = -918.1111000000001
$ws.Range("H77").Value = 1829.8182
$ws.Range("I77").Value = 1792.1111
$ws.Range("K77").Value = 8960.5555
$ws.Range("M77").Value = -4592.5555
$ws.Range("H122").Value = 5012.25
$ws.Range("J122").Value = 5012
$ws.Range("L122").Value = 15036
$ws.Range("N122").Value = -19936
$ws.Range("H132").Value = 1847
$ws.Range("I132").Value = 1339.1666
$ws.Range("K132").Value = 4017.4998
$ws.Range("M132").Value = -1487.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 8999
$ws.Range("J22").Value = 8999
$ws.Range("L22").Value = 8999
$ws.Range("N22").Value = -9345
$ws.Range("H29").Value = 5000
$ws.Range("I29").Value = 5000
$ws.Range("K29").Value = 5000
$ws.Range("M29").Value = -4711
$ws.Range("H36").Value = 8017
$ws.Range("I36").Value = 8017
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 8017
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -7483
$ws.Range("N36").ClearContents()
$ws.Range("H86").Value = 3450
$ws.Range("I86").Value = 3600
$ws.Range("K86").Value = 3600
$ws.Range("M86").Value = -2477
$ws.Range("H89").Value = 3450
$ws.Range("I89").Value = 3600
$ws.Range("K89").Value = 18000
$ws.Range("M89").Value = -12384
$ws.Range("H134").Value = 3133.842
$ws.Range("I134").Value = 3043.3845
$ws.Range("K134").Value = 9130.1535
$ws.Range("M134").Value = -6595.1535

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1164.8334
$ws.Range("I16").Value = 997
$ws.Range("K16").Value = 997
$ws.Range("M16").Value = -710
$ws.Range("H31").Value = 3124.0527
$ws.Range("I31").Value = 1201.4286
$ws.Range("J31").Value = 4245.5835
$ws.Range("K31").Value = 1201.4286
$ws.Range("L31").Value = 4245.5835
$ws.Range("M31").Value = -906.4286
$ws.Range("N31").Value = -4835.5835
$ws.Range("H34").Value = 3124.0527
$ws.Range("I34").Value = 1201.4286
$ws.Range("J34").Value = 4245.5835
$ws.Range("K34").Value = 1201.4286
$ws.Range("L34").Value = 4245.5835
$ws.Range("M34").Value = -999.4286
$ws.Range("N34").Value = -4649.5835
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H60").Value = 20000
$ws.Range("I60").Value = 20000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 20000
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -19489
$ws.Range("N60").ClearContents()
$ws.Range("H62").Value = 13700
$ws.Range("J62").Value = 2000
$ws.Range("L62").Value = 2000
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 13700
$ws.Range("J65").Value = 2000
$ws.Range("L65").Value = 10000
$ws.Range("N65").Value = -16240
$ws.Range("H99").Value = 48077
$ws.Range("I99").Value = 7405
$ws.Range("K99").Value = 7405
$ws.Range("M99").Value = -5907
$ws.Range("H103").Value = 5464.8335
$ws.Range("I103").Value = 6088.8
$ws.Range("K103").Value = 6088.8
$ws.Range("M103").Value = -4916.8
$ws.Range("H113").Value = 1164.8334
$ws.Range("I113").Value = 997
$ws.Range("K113").Value = 997
$ws.Range("M113").Value = 1173
$ws.Range("H126").Value = 48077
$ws.Range("I126").Value = 7405
$ws.Range("K126").Value = 22215
$ws.Range("M126").Value = -19745
$ws.Range("H134").Value = 2828.842
$ws.Range("I134").Value = 2254.9412
$ws.Range("K134").Value = 6764.823600000001
$ws.Range("M134").Value = -4229.823600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 600
$ws.Range("M2").Value = -487
$ws.Range("H12").Value = 1160
$ws.Range("I12").Value = 74.833336
$ws.Range("J12").Value = 1566.9375
$ws.Range("K12").Value = 224.500008
$ws.Range("L12").Value = 4700.8125
$ws.Range("M12").Value = -51.50000800000001
$ws.Range("N12").Value = -5046.8125
$ws.Range("H51").Value = 1443.4
$ws.Range("I51").Value = 1306
$ws.Range("J51").Value = 1993
$ws.Range("K51").Value = 3918
$ws.Range("L51").Value = 5979
$ws.Range("M51").Value = -3458
$ws.Range("N51").Value = -6899
$ws.Range("H107").Value = 1223.08
$ws.Range("J107").Value = 1281.1765
$ws.Range("L107").Value = 3843.5295
$ws.Range("N107").Value = -7683.529500000001
$ws.Range("H109").Value = 1999
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("H129").Value = 4136.067
$ws.Range("I129").Value = 862.2
$ws.Range("J129").Value = 5773
$ws.Range("K129").Value = 2586.6
$ws.Range("L129").Value = 17319
$ws.Range("M129").Value = 2413.4
$ws.Range("N129").Value = -27319
$ws.Range("H131").Value = 3017349.5
$ws.Range("J131").Value = 3473678.5
$ws.Range("L131").Value = 10421035.5
$ws.Range("N131").Value = -10431115.5
$ws.Range("H133").Value = 7625
$ws.Range("I133").Value = 4000
$ws.Range("J133").Value = 8142.857
$ws.Range("K133").Value = 12000
$ws.Range("L133").Value = 24428.571
$ws.Range("M133").Value = -6940
$ws.Range("N133").Value = -34548.571
$ws.Range("H134").Value = 6499.3
$ws.Range("I134").Value = 2999.6
$ws.Range("K134").Value = 8998.799999999999
$ws.Range("M134").Value = -3928.799999999999
$ws.Range("H139").Value = 9938.556
$ws.Range("I139").Value = 6944.8335
$ws.Range("J139").Value = 15926
$ws.Range("K139").Value = 20834.5005
$ws.Range("L139").Value = 47778
$ws.Range("M139").Value = -15694.5005
$ws.Range("N139").Value = -58058
$ws.Range("H140").Value = 3150.4285
$ws.Range("I140").Value = 3150.4285
$ws.Range("K140").Value = 9451.2855
$ws.Range("M140").Value = -4271.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2848.5
$ws.Range("I113").Value = 2557.8
$ws.Range("J113").Value = 3333
$ws.Range("K113").Value = 2557.8
$ws.Range("L113").Value = 3333
$ws.Range("M113").Value = -387.8000000000002
$ws.Range("N113").Value = -7673
$ws.Range("H122").Value = 2528.25
$ws.Range("I122").Value = 2760.3333
$ws.Range("J122").Value = 2389
$ws.Range("K122").Value = 8280.999899999999
$ws.Range("L122").Value = 7167
$ws.Range("M122").Value = -5830.999899999999
$ws.Range("N122").Value = -12067

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3661.889
$ws.Range("I7").Value = 2513.9
$ws.Range("J7").Value = 5096.875
$ws.Range("K7").Value = 2513.9
$ws.Range("L7").Value = 5096.875
$ws.Range("M7").Value = -2401.9
$ws.Range("N7").Value = -5320.875
$ws.Range("H53").Value = 33500
$ws.Range("I53").Value = 5000
$ws.Range("J53").Value = 62000
$ws.Range("K53").Value = 5000
$ws.Range("L53").Value = 62000
$ws.Range("M53").Value = -4482
$ws.Range("N53").Value = -63036
$ws.Range("H61").Value = 9416.286
$ws.Range("I61").Value = 9321
$ws.Range("J61").Value = 9487.75
$ws.Range("K61").Value = 9321
$ws.Range("L61").Value = 9487.75
$ws.Range("M61").Value = -9119
$ws.Range("N61").Value = -9891.75
$ws.Range("H97").Value = 10748.6
$ws.Range("J97").Value = 10748.6
$ws.Range("L97").Value = 10748.6
$ws.Range("N97").Value = -12730.6
$ws.Range("H113").Value = 9416.286
$ws.Range("I113").Value = 9321
$ws.Range("J113").Value = 9487.75
$ws.Range("K113").Value = 9321
$ws.Range("L113").Value = 9487.75
$ws.Range("M113").Value = -7151
$ws.Range("N113").Value = -13827.75
$ws.Range("H126").Value = 3661.889
$ws.Range("I126").Value = 2513.9
$ws.Range("J126").Value = 5096.875
$ws.Range("K126").Value = 7541.700000000001
$ws.Range("L126").Value = 15290.625
$ws.Range("M126").Value = -5071.700000000001
$ws.Range("N126").Value = -20230.625
$ws.Range("H132").Value = 3699
$ws.Range("J132").Value = 3784.4285
$ws.Range("L132").Value = 11353.2855
$ws.Range("N132").Value = -16413.2855
$ws.Range("H136").Value = 8176.8096
$ws.Range("I136").Value = 5653.8887
$ws.Range("J136").Value = 10069
$ws.Range("K136").Value = 16961.6661
$ws.Range("L136").Value = 30207
$ws.Range("M136").Value = -14411.6661
$ws.Range("N136").Value = -35307

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 40012
$ws.Range("I40").Value = 40012
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 40012
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -39863
$ws.Range("N40").ClearContents()
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("H46").Value = 46749.5
$ws.Range("I46").Value = 25999
$ws.Range("K46").Value = 25999
$ws.Range("M46").Value = -25768
$ws.Range("H107").Value = 710.5
$ws.Range("I107").Value = 747.86664
$ws.Range("J107").Value = 150
$ws.Range("K107").Value = 2243.59992
$ws.Range("L107").Value = 450
$ws.Range("M107").Value = -323.5999199999997
$ws.Range("N107").Value = -4290
$ws.Range("H113").Value = 859.1053000000001
$ws.Range("I113").Value = 877.9091
$ws.Range("J113").Value = 833.25
$ws.Range("K113").Value = 2633.7273
$ws.Range("L113").Value = 2499.75
$ws.Range("M113").Value = -463.7273
$ws.Range("N113").Value = -6839.75
$ws.Range("H126").Value = 2641.383
$ws.Range("I126").Value = 1888.6154
$ws.Range("J126").Value = 6311.125
$ws.Range("K126").Value = 5665.8462
$ws.Range("L126").Value = 18933.375
$ws.Range("M126").Value = -3195.8462
$ws.Range("N126").Value = -23873.375
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 5785.552
$ws.Range("I132").Value = 4115.6665
$ws.Range("J132").Value = 13801
$ws.Range("K132").Value = 12346.9995
$ws.Range("L132").Value = 41403
$ws.Range("M132").Value = -9816.999500000002
$ws.Range("N132").Value = -46463
$ws.Range("H134").Value = 46749.5
$ws.Range("I134").Value = 25999
$ws.Range("K134").Value = 77997
$ws.Range("M134").Value = -75462
$ws.Range("H135").Value = 64107.5
$ws.Range("J135").Value = 64107.5
$ws.Range("L135").Value = 64107.5
$ws.Range("N135").Value = -74247.5
$ws.Range("H136").Value = 3182.7144
$ws.Range("I136").Value = 2366.6667
$ws.Range("K136").Value = 7100.000100000001
$ws.Range("M136").Value = -4550.000100000001
